$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.786
$ws.Range("C2").Value = 0.877
$ws.Range("D2").Value = 0.945
$ws.Range("E2").Value = 0.83
$ws.Range("F2").Value = 0.44
$ws.Range("G2").Value = 0.8179999999999999

$ws.Range("B3").Value = 0.72
$ws.Range("C3").Value = 0.779
$ws.Range("D3").Value = 0.909
$ws.Range("E3").Value = 0.796
$ws.Range("F3").Value = 0.261
$ws.Range("G3").Value = 0.864

$ws.Range("B4").Value = 0.702
$ws.Range("C4").Value = 0.887
$ws.Range("D4").Value = 0.912
$ws.Range("E4").Value = 0.7
$ws.Range("F4").Value = 0.379
$ws.Range("G4").Value = 0.779
